$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.830.71'
$ws.Range("E2").Value = '  +1.53%  '
$ws.Range("D3").Value = '3.886.95'
$ws.Range("E3").Value = '  +3.27%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''466.27'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +9.41%  '
$ws.Range("D6").Value = '''142.49'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.66%  '
$ws.Range("D7").Value = '''0.622'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '''0.732'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("E10").Value = '  +8.82%  '
$ws.Range("E11").Value = '  +9.98%  '
$ws.Range("D12").Value = '''42.79'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.45%  '
$ws.Range("D13").Value = '4.506.21'
$ws.Range("E13").Value = '  +3.42%  '
$ws.Range("E14").Value = '  -0.86%  '
$ws.Range("D15").Value = '''14.88'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.99%  '
$ws.Range("D16").Value = '3.871.93'
$ws.Range("E16").Value = '  +3.31%  '
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").Value = '''19.74'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("E19").Value = '  +4.01%  '
$ws.Range("D20").Value = '67.044.80'
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("D21").Value = '''427.65'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +6.65%  '
$ws.Range("D22").Value = '''14.68'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.13%  '
$ws.Range("D23").Value = '''3.32'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.88%  '
$ws.Range("D24").Value = '''87.95'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.68%  '
$ws.Range("B25").Value = 'EthereumClassic'
$ws.Range("C25").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D25").Value = '''38.32'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +5.04%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").Value = '''3.53'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +9.16%  '
$ws.Range("D27").Value = '''5.77'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +6.76%  '
$ws.Range("D28").Value = '''9.96'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.86%  '
$ws.Range("D29").Value = '''9.64'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.55%  '
$ws.Range("D30").Value = '''730.22'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.19%  '
$ws.Range("D31").Value = '''13.70'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.46%  '
$ws.Range("E32").Value = '  -0.37%  '
$ws.Range("D33").Value = '''2.77'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("D34").Value = '''43.19'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +6.02%  '
$ws.Range("E35").Value = '  +5.51%  '
$ws.Range("D36").Value = '''57.30'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.48%  '
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("D38").Value = '''5.38'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.62%  '
$ws.Range("D39").Value = '0.0₃0769'
$ws.Range("E39").Value = '  +16.73%  '
$ws.Range("D40").Value = '''0.0473'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.25%  '
$ws.Range("E41").Value = '  +8.15%  '
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").Value = '''2.58'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.82%  '
$ws.Range("E43").Value = '  +0.64%  '
$ws.Range("D44").Value = '''0.337'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +6.41%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '''1.00'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("E46").Value = '  +5.52%  '
$ws.Range("E47").Value = '  +6.04%  '
$ws.Range("D48").Value = '''3.38'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.72%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '''144.23'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.58%  '
$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").Value = '''3.12'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.91%  '
$ws.Range("D51").Value = '''2.86'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.85%  '
